$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address -> new literal text value.
# NumberFormat "@" forces the assignment to be stored as text (matching the
# original inlineStr cells) instead of being auto-coerced to a number; the
# Style reset afterwards puts the cell back on the workbook's default style
# record so no formatting/style actually changes.
$updates = @{
    "D2" = "65.810.03"
    "E2" = "  +1.46%  "
    "D3" = "2.695.14"
    "E3" = "  +2.14%  "
    "E4" = "  +0.05%  "
    "D5" = "605.25"
    "E5" = "  +2.04%  "
    "D6" = "157.98"
    "E7" = "  +0.02%  "
    "D8" = "0.587"
    "E9" = "  +6.15%  "
    "D10" = "6.05"
    "E10" = "  +4.88%  "
    "D11" = "0.402"
    "E11" = "  +1.20%  "
    "E12" = "  +1.11%  "
    "D13" = "30.17"
    "E13" = "  +4.40%  "
    "E14" = "  +9.63%  "
    "D15" = "3.180.12"
    "E15" = "  +2.20%  "
    "D16" = "65.665.92"
    "E16" = "  +1.48%  "
    "D17" = "2.702.97"
    "E17" = "  +1.25%  "
    "E18" = "  +1.04%  "
    "D19" = "4.88"
    "D20" = "359.50"
    "E20" = "  +2.49%  "
    "D21" = "7.52"
    "E21" = "  +3.89%  "
    "E22" = "  -0.13%  "
    "E23" = "  +3.87%  "
    "D24" = "9.82"
    "E24" = "  +3.97%  "
    "E25" = "  +13.53%  "
    "E26" = "  -3.20%  "
    "E27" = "  +3.43%  "
    "E28" = "  +5.55%  "
    "D29" = "8.38"
    "E29" = "  +3.54%  "
    "D30" = "544.20"
    "E30" = "  +6.54%  "
    "E31" = "  +4.66%  "
    "E32" = "  +0.09%  "
    "E33" = "  +1.06%  "
    "D34" = "6.72"
    "E34" = "  +8.03%  "
    "E35" = "  -3.64%  "
    "D36" = "0.433"
    "E36" = "  +2.35%  "
    "D37" = "20.79"
    "E37" = "  +3.49%  "
    "D38" = "162.94"
    "E38" = "  -1.21%  "
    "D39" = "1.99"
    "E39" = "  -0.70%  "
    "D40" = "1.00"
    "E40" = "  +0.02%  "
    "D42" = "42.90"
    "E42" = "  +1.55%  "
    "D43" = "168.94"
    "E43" = "  +2.85%  "
    "E44" = "  +3.00%  "
    "E45" = "  +0.70%  "
    "E46" = "  +3.66%  "
    "D47" = "2.28"
    "E47" = "  +4.29%  "
    "E48" = "  +5.15%  "
    "D49" = "0.659"
    "E49" = "  +2.05%  "
    "D50" = "21.07"
    "E50" = "  +9.34%  "
    "D51" = "0.0993"
    "E51" = "  +1.43%  "
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
